# ---------------------------------------------------------------------------
# lab_06/data/lab_data.xlsx -- "Add day 2 files"
#
#   * one_tone         : add an "Attenuation [dB]" column (C) = 30 everywhere
#   * two_tone         : insert an "Attenuation [dB]" column (B) and weave in
#                        the new "day 2" rows (10 dB sweeps) between the
#                        original (30 dB) rows
#   * two_tone_double  : brand-new sheet with the double-tone/attenuation data
#   * assorted view cosmetics (zoom, selection, active sheet, column widths)
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ===========================================================================
# 1. one_tone -- new "Attenuation [dB]" column
# ===========================================================================
$oneTone = $wb.Worksheets.Item("one_tone")
$oneTone.Range("C1").Value = "Attenuation [dB]"
for ($r = 2; $r -le 16; $r++) {
    $oneTone.Cells.Item($r, 3).Value = 30
}

$oneTone.Activate()
$oneTone.Range("D14").Select()
$excel.ActiveWindow.Zoom = 214

# ===========================================================================
# 2. two_tone -- insert "Attenuation [dB]" column + weave in new rows
# ===========================================================================
$twoTone = $wb.Worksheets.Item("two_tone")

# Shift the tone-pair data from B:J -> C:K, freeing up column B.
$twoTone.Columns.Item(2).Insert()
$twoTone.Range("B1").Value = "Attenuation [dB]"

# The 30 dB sweeps already present only had blanks in the (new) B column --
# fill them in now while the row numbers still match the original rows.
$twoTone.Cells.Item(2, 2).Value = 30
$twoTone.Cells.Item(3, 2).Value = 30
$twoTone.Cells.Item(4, 2).Value = 30
$twoTone.Cells.Item(5, 2).Value = 30
$twoTone.Cells.Item(6, 2).Value = 30
$twoTone.Cells.Item(7, 2).Value = 30
$twoTone.Cells.Item(8, 2).Value = 30
$twoTone.Cells.Item(9, 2).Value = 30
$twoTone.Cells.Item(10, 2).Value = 30

# Weave in the new attenuation=10dB rows -- one freshly-inserted row ahead of
# each of the (then-current) rows 6, 8, 10 and 12.
$twoTone.Rows.Item(6).Insert()
$twoTone.Rows.Item(8).Insert()
$twoTone.Rows.Item(10).Insert()
$twoTone.Rows.Item(12).Insert()

# Drop the stray trailing blank row (the old A11 placeholder, now pushed to
# row 15).
$twoTone.Rows.Item(15).Delete()

function Set-Row($row, $values) {
    foreach ($colLetter in $values.Keys) {
        $twoTone.Range("$colLetter$row").Value = $values[$colLetter]
    }
}

Set-Row 6  @{ A = -25; B = 10; H = -73.7;               I = -63.63; J = -63.6;  K = -73.400000000000006 }
Set-Row 8  @{ A = -22; B = 10; H = -64.099999999999994; I = -54.6;  J = -54.5;  K = -63.59 }
Set-Row 10 @{ A = -20; B = 10; H = -57.5;                I = -48.5;  J = -48.3;  K = -56.9 }
Set-Row 12 @{ A = -19; B = 10; H = -53;                  I = -45;    J = -44;    K = -53 }

# Residual formatted-but-empty placeholder cells that show up at H21:I22 in
# the saved file (carried-over selection/format artifact).
$fmtSrc = $twoTone.Range("A2")
$fmtSrc.Copy()
$twoTone.Range("H21:I22").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Column widths for A/B (custom-widened to fit the new attenuation column).
$twoTone.Columns.Item(1).ColumnWidth = 14.83
$twoTone.Columns.Item(2).ColumnWidth = 12

$twoTone.Activate()
$twoTone.Range("K1").Select()
$excel.ActiveWindow.Zoom = 133

# ===========================================================================
# 3. two_tone_double -- brand-new sheet
# ===========================================================================
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$dbl = $wb.Worksheets.Add($null, $lastSheet)
$dbl.Name = "two_tone_double"

$dbl.Range("A1").Value = "Input Power [dBm]"
$dbl.Range("B1").Value = "Attenuation [dB]"
$dbl.Cells.Item(1, 3).Value = 600
$dbl.Cells.Item(1, 4).Value = 601
$dbl.Cells.Item(1, 5).Value = 1800
$dbl.Cells.Item(1, 6).Value = 1801
$dbl.Cells.Item(1, 7).Value = 1802
$dbl.Cells.Item(1, 8).Value = 1803

function Set-DblRow($row, $a, $b, $c, $d, $e, $f, $g, $h) {
    $dbl.Cells.Item($row, 1).Value = $a
    $dbl.Cells.Item($row, 2).Value = $b
    $dbl.Cells.Item($row, 3).Value = $c
    $dbl.Cells.Item($row, 4).Value = $d
    $dbl.Cells.Item($row, 5).Value = $e
    $dbl.Cells.Item($row, 6).Value = $f
    $dbl.Cells.Item($row, 7).Value = $g
    $dbl.Cells.Item($row, 8).Value = $h
}

Set-DblRow 2 -50   10 -16                 -15.91               -44.1 -31.7               -31.6 -43.5
Set-DblRow 3 -53   10 -16.95               -16.600000000000001 -46.12 -38.200000000000003 -38.1 -45.9
Set-DblRow 4 -55   10 -18.079999999999998  -18.010000000000002 -53.3 -45.8               -46.8 -53.1
Set-DblRow 5 -57   10 -19.66               -19.57               -61.4 -53.6               -53.9 -62.2
Set-DblRow 6 -58.2 10 -20.74               -20.66               -66.569999999999993 -58.3 -58.5 -66.5

$dbl.Columns.Item(1).ColumnWidth = 14.5
$dbl.Columns.Item(2).ColumnWidth = 13.17

$dbl.Activate()
$dbl.Range("E3:H3").Select()
$excel.ActiveWindow.Zoom = 134

# ===========================================================================
# 4. sweep_time -- view-only cosmetics (zoom; no longer the active tab)
# ===========================================================================
$sweepTime = $wb.Worksheets.Item("sweep_time")
$sweepTime.Activate()
$excel.ActiveWindow.Zoom = 182

# Leave two_tone_double as the final active sheet/tab.
$dbl.Activate()

Write-Host "done"
